$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    3 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 3.536033448013082)
    4 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 6.82939032824165)
    5 = @(0.0006075818656279264, 1.626987699542094, 18.71679738969934, 13.86384647080068, 0, 34.20823914190774)
    6 = @(1.445647641019636, 9.983522426115931, 189.6080260415259, 2797.565817734744, 0, 2998.603013843405)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
